$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(7, 2).Value = "col1"
$ws.Cells.Item(7, 3).Value = "col2"
$ws.Cells.Item(7, 4).Value = "col3"

# Data row 1
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 3).Value = "test_col2"
$ws.Cells.Item(8, 4).Value = "test_col3"

# Data row 2
$ws.Cells.Item(9, 2).Value = 2
$ws.Cells.Item(9, 3).Value = "test_col21"
$ws.Cells.Item(9, 4).Value = "test_col31"

# Set the active cell/selection to D8, matching the target workbook state
$ws.Range("D8").Select()
